# Regenerate orders with updated distance/sizes
# Mapping:
#   Distance: D64 -> D69, D80 -> D86, D51 -> D55
#   Size:     S30 -> S31 (S25, S20 unchanged)
# Applies to columns: Condition (B), Filename_Left (D), Filename_Right (E),
#                      Distance (H), Size (J)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rng = $ws.UsedRange
$lastRow = $rng.Rows.Count

function Convert-Token([string]$value) {
    if ($null -eq $value) { return $value }
    $out = $value
    $out = $out.Replace("D64", "D69")
    $out = $out.Replace("D80", "D86")
    $out = $out.Replace("D51", "D55")
    $out = $out.Replace("S30", "S31")
    return $out
}

$cols = @(2, 4, 5, 8, 10)

for ($r = 2; $r -le $lastRow; $r++) {
    foreach ($c in $cols) {
        $cell = $ws.Cells.Item($r, $c)
        $val = $cell.Value()
        if ($val -is [string]) {
            $newVal = Convert-Token $val
            if ($newVal -ne $val) {
                $cell.Value = $newVal
            }
        }
    }
}
